# fix(scoring): update PSS labels to match Excel format
#
# The PSS sheet's row 2 (A2:R2) held a shared formula
#   =IF(A1>26,"Stress Berat",IF(A1>13,"Stress Sedang","Stress Ringan"))
# which (because row 1 holds text, not the numeric scores) always
# resolved to "Stress Berat". Replace each cell with the correct static
# label ("Stress Ringan" / "Stress Sedang") so the sheet matches the
# validated scoring output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PSS")
$ws.Activate()

$labels = @{
    "A2" = "Stress Ringan"
    "B2" = "Stress Ringan"
    "C2" = "Stress Ringan"
    "D2" = "Stress Ringan"
    "E2" = "Stress Ringan"
    "F2" = "Stress Ringan"
    "G2" = "Stress Ringan"
    "H2" = "Stress Ringan"
    "I2" = "Stress Ringan"
    "J2" = "Stress Ringan"
    "K2" = "Stress Ringan"
    "L2" = "Stress Ringan"
    "M2" = "Stress Ringan"
    "N2" = "Stress Sedang"
    "O2" = "Stress Ringan"
    "P2" = "Stress Ringan"
    "Q2" = "Stress Ringan"
    "R2" = "Stress Sedang"
}

foreach ($addr in $labels.Keys) {
    $ws.Range($addr).Value = $labels[$addr]
}

# Mirror the author's final selection/scroll state on the sheet.
$ws.Range("Q12").Select()

Write-Host "PSS row 2 labels updated to static Stress Ringan/Stress Sedang values"
